$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddProductCategory1")
$ws.Activate()

# Update row 1 height (16.25 -> 15.65)
$ws.Rows.Item(1).RowHeight = 15.65

# New "Company" column J (values 20..28), and new "Product" values in column I (13..18 for rows 4-9)
$jValues = @(20, 21, 22, 23, 24, 25, 26, 27, 28)
for ($i = 0; $i -lt $jValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}

$iValues = @(13, 14, 15, 16, 17, 18)
for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 4
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
}

# Normalize date format for G4:G9 so they share the same style as G1:G3 (drops duplicate DD/MM/YY numFmt)
$ws.Range("G4:G9").NumberFormat = "DD/MM/YY"

# Move the active selection to I4 (was I19)
$ws.Range("I4").Select()
